$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values like "31.182.90" or
# "0.05000" are not coerced into numbers, losing formatting/precision.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '31.182.90'
$ws.Range("E2").Value = '  +1.88%  '
$ws.Range("D3").Value = '1.959.08'
$ws.Range("E3").Value = '  +1.05%  '
$ws.Range("E4").Value = '  +0.31%  '
$ws.Range("D5").Value = '246.93'
$ws.Range("E5").Value = '  +0.19%  '
$ws.Range("E6").Value = '  +0.25%  '
$ws.Range("D7").Value = '0.4894'
$ws.Range("E7").Value = '  +1.06%  '
$ws.Range("D8").Value = '0.2971'
$ws.Range("E8").Value = '  +1.67%  '
$ws.Range("D9").Value = '0.06839'
$ws.Range("E9").Value = '  +0.16%  '
$ws.Range("D10").Value = '19.09'
$ws.Range("E10").Value = '  -1.92%  '
$ws.Range("D11").Value = '106.56'
$ws.Range("E11").Value = '  -5.20%  '
$ws.Range("D12").Value = '1.937.47'
$ws.Range("E12").Value = '  +0.02%  '
$ws.Range("D13").Value = '0.07750'
$ws.Range("E13").Value = '  +2.12%  '
$ws.Range("D14").Value = '5.418'
$ws.Range("E14").Value = '  -1.20%  '
$ws.Range("D15").Value = '0.7159'
$ws.Range("E15").Value = '  +4.86%  '
$ws.Range("D16").Value = '284.77'
$ws.Range("E16").Value = '  -4.67%  '
$ws.Range("D17").Value = '31.049.09'
$ws.Range("E17").Value = '  +1.52%  '
$ws.Range("D18").Value = '0.000007760'
$ws.Range("E18").Value = '  +1.12%  '
$ws.Range("D19").Value = '13.23'
$ws.Range("E19").Value = '  +0.79%  '
$ws.Range("B20").Value = 'Dai'
$ws.Range("C20").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D20").Value = '1.002'
$ws.Range("E20").Value = '  +0.30%  '
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '2.191.77'
$ws.Range("E21").Value = '  +0.31%  '
$ws.Range("D22").Value = '5.525'
$ws.Range("E22").Value = '  -1.34%  '
$ws.Range("D23").Value = '1.002'
$ws.Range("E23").Value = '  +0.29%  '
$ws.Range("D24").Value = '6.610'
$ws.Range("E24").Value = '  +1.30%  '
$ws.Range("D25").Value = '9.922'
$ws.Range("E25").Value = '  +4.00%  '
$ws.Range("D26").Value = '168.95'
$ws.Range("E26").Value = '  +0.84%  '
$ws.Range("D27").Value = '19.93'
$ws.Range("E27").Value = '  -2.48%  '
$ws.Range("D28").Value = '2.209'
$ws.Range("E28").Value = '  +2.81%  '
$ws.Range("E29").Value = '  +0.32%  '
$ws.Range("D30").Value = '0.1052'
$ws.Range("E30").Value = '  -1.83%  '
$ws.Range("D31").Value = '4.737'
$ws.Range("E31").Value = '  +15.24%  '
$ws.Range("D32").Value = '4.499'
$ws.Range("E32").Value = '  +7.86%  '
$ws.Range("D33").Value = '0.05000'
$ws.Range("E33").Value = '  -0.34%  '
$ws.Range("D34").Value = '0.7625'
$ws.Range("E34").Value = '  +2.08%  '
$ws.Range("D35").Value = '1.169'
$ws.Range("E35").Value = '  +1.47%  '
$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D36").Value = '2.736'
$ws.Range("E36").Value = '  +0.80%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = '0.02048'
$ws.Range("E37").Value = '  +0.14%  '
$ws.Range("D38").Value = '2.715'
$ws.Range("E38").Value = '  +0.66%  '
$ws.Range("D39").Value = '2.160'
$ws.Range("E39").Value = '  +5.91%  '
$ws.Range("D40").Value = '6.448'
$ws.Range("E40").Value = '  +10.16%  '
$ws.Range("D41").Value = '0.8846'
$ws.Range("E41").Value = '  +1.05%  '
$ws.Range("D42").Value = '109.79'
$ws.Range("E42").Value = '  -0.36%  '
$ws.Range("D43").Value = '0.4465'
$ws.Range("E43").Value = '  -0.28%  '
$ws.Range("D44").Value = '72.82'
$ws.Range("E44").Value = '  +4.01%  '
$ws.Range("D45").Value = '1.001'
$ws.Range("E45").Value = '  +0.01%  '
$ws.Range("D46").Value = '7.512'
$ws.Range("E46").Value = '  +2.57%  '
$ws.Range("D47").Value = '988.98'
$ws.Range("E47").Value = '  +16.83%  '
$ws.Range("D48").Value = '0.1275'
$ws.Range("E48").Value = '  +3.01%  '
$ws.Range("D49").Value = '9.369'
$ws.Range("E49").Value = '  +0.14%  '
$ws.Range("D50").Value = '0.2610'
$ws.Range("E50").Value = '  +2.39%  '
$ws.Range("D51").Value = '35.97'
$ws.Range("E51").Value = '  +2.36%  '

# Restore default style on column D now that text values are safely stored,
# so no stray style attributes linger on the cells.
$priceRange.Style = "Normal"

